# Applies the "Penalty Reward System" edit to the PO data workbook.
#
# Changes made:
#  - On the "Weekly Quantity" sheet:
#      * The requested quantity for the week of 45130.99999999999 is
#        corrected from 140 to 100.
#      * Three weekly rows (45088.99999999999, 45102.99999999999 and
#        45137.99999999999) are removed entirely, shifting the rows
#        below them upward.
#  - On the "Monthly Trend" sheet, the two affected monthly totals are
#    updated to reflect the corrected/removed weekly data:
#      * 45107.99999999999 total drops from 640 to 420.
#      * 45138.99999999999 total drops from 440 to 220.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# 1) Fix the requested quantity value for the 45130.99999999999 week
#    (row 9) while row numbers are still in their original positions.
$wsWeekly.Cells.Item(9, 2).Value = 100

# 2) Remove the three obsolete weekly rows. Delete from the bottom up
#    so the row numbers of the remaining deletions stay correct.
$wsWeekly.Rows.Item(10).Delete()
$wsWeekly.Rows.Item(5).Delete()
$wsWeekly.Rows.Item(3).Delete()

# 3) Update the corresponding monthly totals.
$wsMonthly.Cells.Item(3, 2).Value = 420
$wsMonthly.Cells.Item(4, 2).Value = 220
